$wb = $excel.ActiveWorkbook

# --- "Special Character Removal" sheet: add the Orange ("J:L") block ---
$ws = $wb.Worksheets.Item("Special Character Removal")

# Header row (row 1) gets an extra blank, wrap-formatted cell under the new column L
$ws.Range("L1").Value2 = ""

# Row 2 - column headers
$ws.Range("J2").Value2 = "Action"
$ws.Range("K2").Value2 = "Time"
$ws.Range("L2").Value2 = "Content"

# Row 3
$ws.Range("J3").Value2 = "Load Data"
$ws.Range("K3").Value2 = "1 min"
$ws.Range("L3").Value2 = "Upload the dataset to the 'Import' page."

# Row 4
$ws.Range("J4").Value2 = "Inspect Data"
$ws.Range("K4").Value2 = "2 min"
$ws.Range("L4").Value2 = "Spot special characters on the 'Data Profiling' page."

# Row 5
$ws.Range("J5").Value2 = "Preprocess"
$ws.Range("K5").Value2 = "1 min"
$ws.Range("L5").Value2 = "Remove special characters on the 'Data Cleaning' page."

# Row 6
$ws.Range("J6").Value2 = "Verify Changes"
$ws.Range("K6").Value2 = "1 min"
$ws.Range("L6").Value2 = "Inspect cleaned data on the 'Data Profiling' page."

# Row 7 - overall total
$ws.Range("J7").Value2 = "Overall"
$ws.Range("K7").Value2 = "5 min"

# Match styles used by the analogous J:L block on the "Rename Column" sheet
$ws.Range("J2,K2").Style = "Normal"
$ws.Range("J2:L2").Font.Bold = $true
$ws.Range("J3:K6").WrapText = $true
$ws.Range("L3:L6").WrapText = $true
$ws.Range("J7:K7").Font.Bold = $true

# Row heights to match the new wrapped content
$ws.Rows.Item(3).RowHeight = 54
$ws.Rows.Item(4).RowHeight = 72
$ws.Rows.Item(5).RowHeight = 72
$ws.Rows.Item(6).RowHeight = 72

# Column widths for the two newly-used columns
$ws.Columns.Item(10).ColumnWidth = 13.330729166666666
$ws.Columns.Item(12).ColumnWidth = 14.166666666666666

# Make this sheet the active / selected tab, with L1:L1048576 selected
$ws.Activate()
$ws.Range("L1:L1048576").Select()

# --- "Rename Column" sheet loses the tabSelected flag (handled automatically
#     by activating "Special Character Removal" above, which Excel treats as
#     mutually exclusive single-tab selection) ---
